$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right before the current row 69 ("Vega Modelo de
# Temuco" / Acelga daily series). This shifts every existing row from 69
# downward to 70.. and grows the sheet from A1:R161 to A1:R162.
$ws.Rows("69:69").Insert()

# Populate the newly inserted row 69 with the new daily observation.
$ws.Cells.Item(69, 1).Value = 10
$ws.Cells.Item(69, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(69, 3).Value = "La Araucanía"
$ws.Cells.Item(69, 4).Value = 44413
$ws.Cells.Item(69, 5).Value = 9
$ws.Cells.Item(69, 6).Value = 100112009
$ws.Cells.Item(69, 7).Value = "Acelga"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 40
$ws.Cells.Item(69, 11).Value = 12000
$ws.Cells.Item(69, 12).Value = 12000
$ws.Cells.Item(69, 13).Value = 12000
$ws.Cells.Item(69, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(69, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(69, 16).Value = 1000
$ws.Cells.Item(69, 17).Value = 12
$ws.Cells.Item(69, 18).Value = "Hortaliza"
